# target invest refinements 2 and preparing data for NL case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# End Year: 2030 -> 2050
$ws.Range("B4").Value = 2050

# maximum_investment_capacity_per_year: 300 -> 1000000
$ws.Range("B13").Value = 1000000

# targetinvestment_per_year: TRUE -> FALSE (preparing data for NL case)
$ws.Range("B20").Value = $false

# Remove the custom font/alignment styling that used to highlight A20,
# reverting the cell back to the workbook's default "Normal" style
$ws.Range("A20").Style = "Normal"

# Bring the sheet into view and leave the selection on the last row
$ws.Activate()
$ws.Range("A20").Select()
